# Add a new "Address" column (F) before the existing "District" column,
# shifting the old F (District) data into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; this shifts the existing District column (F) to G.
$ws.Columns.Item(6).Insert()

# Header for the new column.
$ws.Range("F2").Value = "Address"

# Address values parsed from the original "Name / School, Town, District"
# block that is split across columns B/E/F. Rows whose original text had no
# separate school-address segment (only a District) are left blank.
$ws.Range("F3").Value = "Jeevanprakash H S Old Jewergi Road"
$ws.Range("F6").Value = "Govt. High School Pattan"
$ws.Range("F7").Value = "G H S MadanaSedam"
$ws.Range("F8").Value = "G H S MannurAfzalpur"
$ws.Range("F10").Value = "Milind High SchoolVidya Nagar"
$ws.Range("F11").Value = "Kudal Sangam High SchoolJewargi"
$ws.Range("F12").Value = "G H S HosurAfzalpur"
$ws.Range("F14").Value = "G H S Station GanagapurAfzalpur"
$ws.Range("F15").Value = "Govt. High School AinapurChincholi"
$ws.Range("F16").Value = "G H S Khanadal"
$ws.Range("F17").Value = "G U H S KodaliChincholli"
$ws.Range("F19").Value = "G H S Afzalapur"
$ws.Range("F20").Value = "Govt. H S Bankalaga"
$ws.Range("F23").Value = "Sri. J V V S High School BalurigiAfzalpur"
$ws.Range("F24").Value = "G H S DonagaonChittapur"
$ws.Range("F25").Value = "R G H S ChigaralliJewargi"
$ws.Range("F26").Value = "G H S NagardalaiChincholi"
$ws.Range("F27").Value = "G H S Adaki Sedam"
$ws.Range("F28").Value = "G G H S Sedam"
$ws.Range("F29").Value = "Indiragandhi M High School Brahamapur"
$ws.Range("F30").Value = "N V Vidyaniketan High School Venkatesh Nagar"
$ws.Range("F31").Value = "G H S Afzalpur"
$ws.Range("F32").Value = "G H S KolkundaSedam"
$ws.Range("F33").Value = "Govt. High SchoolMedakSedam"
$ws.Range("F35").Value = "Govt. High School Farhatabad"
$ws.Range("F36").Value = "G H S Jevargi Colony"
$ws.Range("F37").Value = "Govt. High SchoolSindgi"
$ws.Range("F38").Value = "G G H S Jagat"
$ws.Range("F41").Value = "G H S Boys Aland"
$ws.Range("F43").Value = "M S S P HSSirnoor"
$ws.Range("F44").Value = "Mahadevi Girls High School Sharannagar"
$ws.Range("F45").Value = "Sri S S High SchoolSavalagi"
$ws.Range("F46").Value = "Sri Chowdeshwar High SchoolBrahampur"
$ws.Range("F47").Value = "Sri. Sevanikltan High School"
$ws.Range("F49").Value = "Govt. H S Gour(B) Afzalpur taluk"
$ws.Range("F50").Value = "G H S Hunasihadgil"
$ws.Range("F51").Value = "G H S ShellagiChittapur"
$ws.Range("F52").Value = "Goutam High School"
$ws.Range("F53").Value = "G H S ManoorAfzalpur"
$ws.Range("F54").Value = "Govt. Girls High SchoolSulepethChincholi"
$ws.Range("F55").Value = "Govt. P U College (Boys) Chincholi"
$ws.Range("F56").Value = "Govt. High School MashalAfzalpur"
$ws.Range("F57").Value = "Ratnasagar Girls H SC I B Colony"
$ws.Range("F58").Value = "G H S Gudur (SA) Jewargi"
